$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "ENCARGADO" (person in charge) column -- fill the data values first so
# the new shared-string entries are appended in the same order as the
# target file (Jorge, Juan Carlos, Freire, Jair, then ENCARGADO).
$ws.Range("C2").Value = "Jorge"
$ws.Range("C3").Value = "Juan Carlos"
$ws.Range("C6").Value = "Freire"
$ws.Range("C7").Value = "Jair"
$ws.Range("C1").Value = "ENCARGADO"

# Match the header style (row 1) and the data style (rows 2-8, wrap text +
# accent fill + border) already used by column B, instead of synthesizing
# brand-new cell formats.
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

$ws.Range("B2").Copy()
$ws.Range("C2:C8").PasteSpecial(-4122)

# New column width for column C (closest value the engine's pixel-quantized
# ColumnWidth storage can represent near the target 12.28515625 XML width)
$ws.Columns.Item(3).ColumnWidth = 11.5

# Update selection to match the new active cell
$ws.Range("D6").Select()
